$wb = $excel.ActiveWorkbook

# --- Sheet "İş Takip Listesi": shift İşe Başlama / İhale Bitiş dates back by one day ---
$ws1 = $wb.Worksheets.Item("İş Takip Listesi")

for ($r = 2; $r -le 10; $r++) {
    $ws1.Cells.Item($r, 10).Value = "'2025-06-14"
    $ws1.Cells.Item($r, 11).Value = "'2025-11-14"
}

for ($r = 33; $r -le 94; $r++) {
    $ws1.Cells.Item($r, 10).Value = "'2025-06-16"
    $ws1.Cells.Item($r, 11).Value = "'2025-11-16"
}

for ($r = 95; $r -le 122; $r++) {
    $ws1.Cells.Item($r, 10).Value = "'2024-04-14"
    $ws1.Cells.Item($r, 11).Value = "'2025-06-08"
}

# Row 67 status updated to "KESİN ASKIDA"
$ws1.Range("L67").Value = "KESİN ASKIDA"

# --- Sheet "Güncelleme": shift the various tracked dates back by one day ---
$ws2 = $wb.Worksheets.Item("Güncelleme")

$ws2.Range("J2").Value = "'2024-07-21"
$ws2.Range("N2").Value = "'2025-03-24"
$ws2.Range("P2").Value = "'2025-06-10"
$ws2.Range("J3").Value = "'2024-10-22"
$ws2.Range("N3").Value = "'2025-07-12"
$ws2.Range("P3").Value = "'2025-10-20"
$ws2.Range("J4").Value = "'2024-08-26"
$ws2.Range("N4").Value = "'2025-02-16"
$ws2.Range("P4").Value = "'2025-05-13"
$ws2.Range("I5").Value = "'2025-02-19"
$ws2.Range("J6").Value = "'2025-10-02"
$ws2.Range("N6").Value = "'2025-06-21"
$ws2.Range("P6").Value = "'2025-12-17"
$ws2.Range("I7").Value = "'2024-10-22"
$ws2.Range("J7").Value = "'2024-10-22"
$ws2.Range("J8").Value = "'2024-10-10"
$ws2.Range("N8").Value = "'2025-03-11"
$ws2.Range("P8").Value = "'2025-04-13"
$ws2.Range("I9").Value = "'2025-06-07"
$ws2.Range("J9").Value = "'2024-11-23"
$ws2.Range("J10").Value = "'2024-09-21"
$ws2.Range("N10").Value = "'2025-07-02"
$ws2.Range("P10").Value = "'2025-10-24"
$ws2.Range("I11").Value = "'2025-03-29"
$ws2.Range("J11").Value = "'2024-11-04"
$ws2.Range("N11").Value = "'2025-07-22"
$ws2.Range("P11").Value = "'2025-12-17"
$ws2.Range("J12").Value = "'2024-10-02"
$ws2.Range("N12").Value = "'2025-06-11"
$ws2.Range("P12").Value = "'2025-10-14"
$ws2.Range("J13").Value = "'2024-11-30"
$ws2.Range("J14").Value = "'2025-09-28"
$ws2.Range("N14").Value = "'2025-10-20"
$ws2.Range("J15").Value = "'2024-12-19"
$ws2.Range("N15").Value = "'2025-07-09"
$ws2.Range("P15").Value = "'2025-10-21"
$ws2.Range("J16").Value = "'2024-08-17"
$ws2.Range("N16").Value = "'2025-01-25"
$ws2.Range("P16").Value = "'2025-04-13"
$ws2.Range("J17").Value = "'2024-09-02"
$ws2.Range("N17").Value = "'2025-10-20"
$ws2.Range("J18").Value = "'2025-02-09"
$ws2.Range("I19").Value = "'2025-03-30"
$ws2.Range("J19").Value = "'2024-12-19"
$ws2.Range("N19").Value = "'2025-07-29"
$ws2.Range("J20").Value = "'2024-11-30"
$ws2.Range("N20").Value = "'2025-12-04"
$ws2.Range("J21").Value = "'2024-09-23"
$ws2.Range("J22").Value = "'2024-09-23"
$ws2.Range("J23").Value = "'2024-12-01"
$ws2.Range("I24").Value = "'2025-05-28"
$ws2.Range("J25").Value = "'2024-10-26"
$ws2.Range("J27").Value = "'2025-01-16"
$ws2.Range("J28").Value = "'2024-11-14"
$ws2.Range("N28").Value = "'2025-10-30"
$ws2.Range("I29").Value = "'2025-02-03"
$ws2.Range("J29").Value = "'2024-12-01"
$ws2.Range("N29").Value = "'2025-10-11"
